$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) for the rows that changed in the
# source data refresh, in both the "展览" sheet and the "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 567
    $ws.Range("F6").Value = 37
    $ws.Range("F8").Value = 483
    $ws.Range("F9").Value = 3595
    $ws.Range("F10").Value = 59
}
